# Apply "Ran code for averaged intensites on spiral schemes":
# adds 3 new averaging schemes (Gaussian-Quadrature re-grouped + 3 new Spiral
# schemes) to the AlphaFiberF averaged-intensities table, growing it from
# A1:M16 to A1:M19, and refreshes the already-present rows with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 19,13

# row 1
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 2
$data[0,4] = 3
$data[0,5] = 4
$data[0,6] = 5
$data[0,7] = 6
$data[0,8] = 7
$data[0,9] = 8
$data[0,10] = 9
$data[0,11] = 10
$data[0,12] = 11

# row 2 - HKL
$data[1,0] = 0
$data[1,1] = 'HKL'
$data[1,2] = '[1, 1, 0]'
$data[1,3] = '[2, 0, 0]'
$data[1,4] = '[2, 1, 1]'
$data[1,5] = '[2, 2, 0]'
$data[1,6] = '[3, 1, 0]'
$data[1,7] = '[2, 2, 2]'
$data[1,8] = '[3, 2, 1]'
$data[1,9] = '[4, 0, 0]'
$data[1,10] = '2Pairs'
$data[1,11] = '4Pairs'
$data[1,12] = 'MaxUnique'

# row 3 - ND Single
$data[2,0] = 1
$data[2,1] = 'ND Single'
$data[2,2] = 0.98
$data[2,3] = 1.03
$data[2,4] = 0.99
$data[2,5] = 0.98
$data[2,6] = 1.02
$data[2,7] = 1
$data[2,8] = 0.98
$data[2,9] = 1.03
$data[2,10] = 1.01
$data[2,11] = 0.995
$data[2,12] = 1

# row 4 - RD Single
$data[3,0] = 2
$data[3,1] = 'RD Single'
$data[3,2] = 1.03
$data[3,3] = 0.8
$data[3,4] = 1.04
$data[3,5] = 1.03
$data[3,6] = 0.9
$data[3,7] = 1.1
$data[3,8] = 1.04
$data[3,9] = 0.8
$data[3,10] = 0.92
$data[3,11] = 0.9750000000000001
$data[3,12] = 0.985

# row 5 - TD Single
$data[4,0] = 3
$data[4,1] = 'TD Single'
$data[4,2] = 1.03
$data[4,3] = 0.88
$data[4,4] = 1.02
$data[4,5] = 1.03
$data[4,6] = 0.94
$data[4,7] = 1.06
$data[4,8] = 1.02
$data[4,9] = 0.88
$data[4,10] = 0.95
$data[4,11] = 0.99
$data[4,12] = 0.9916666666666666

# row 6 - Morris
$data[5,0] = 4
$data[5,1] = 'Morris'
$data[5,2] = 1.02
$data[5,3] = 0.88
$data[5,4] = 1.03
$data[5,5] = 1.02
$data[5,6] = 0.93
$data[5,7] = 1.09
$data[5,8] = 1.02
$data[5,9] = 0.88
$data[5,10] = 0.9550000000000001
$data[5,11] = 0.9874999999999999
$data[5,12] = 0.9950000000000001

# row 7 - Ring Perpendicular to ND
$data[6,0] = 5
$data[6,1] = 'Ring Perpendicular to ND'
$data[6,2] = 0.999041095890411
$data[6,3] = 1.001643835616438
$data[6,4] = 0.9909589041095891
$data[6,5] = 0.999041095890411
$data[6,6] = 1.000684931506849
$data[6,7] = 0.9838356164383562
$data[6,8] = 0.9917808219178083
$data[6,9] = 1.001643835616438
$data[6,10] = 0.9963013698630137
$data[6,11] = 0.9976712328767123
$data[6,12] = 0.9946575342465752

# row 8 - Ring Perpendicular to RD
$data[7,0] = 6
$data[7,1] = 'Ring Perpendicular to RD'
$data[7,2] = 1.010526315789474
$data[7,3] = 0.9236842105263158
$data[7,4] = 1.011578947368421
$data[7,5] = 1.010526315789474
$data[7,6] = 0.9589473684210527
$data[7,7] = 1.040526315789474
$data[7,8] = 1.011052631578947
$data[7,9] = 0.9236842105263158
$data[7,10] = 0.9676315789473684
$data[7,11] = 0.9890789473684211
$data[7,12] = 0.9927192982456141

# row 9 - Ring Perpendicular to TD
$data[8,0] = 7
$data[8,1] = 'Ring Perpendicular to TD'
$data[8,2] = 0.998421052631579
$data[8,3] = 0.9557894736842105
$data[8,4] = 1.006315789473684
$data[8,5] = 0.998421052631579
$data[8,6] = 0.9752631578947368
$data[8,7] = 1.026315789473684
$data[8,8] = 1.005263157894737
$data[8,9] = 0.9557894736842105
$data[8,10] = 0.9810526315789474
$data[8,11] = 0.9897368421052631
$data[8,12] = 0.9945614035087719

# row 10 - Gaussian-Quadrature
$data[9,0] = 8
$data[9,1] = 'Gaussian-Quadrature'
$data[9,2] = 0.9906828493086272
$data[9,3] = 0.9997385364559017
$data[9,4] = 0.991866189525917
$data[9,5] = 0.9906828493086272
$data[9,6] = 1.001583232288522
$data[9,7] = 0.9882714118751865
$data[9,8] = 0.9929119283421032
$data[9,9] = 0.9997385364559017
$data[9,10] = 0.9958023629909094
$data[9,11] = 0.9932426061497682
$data[9,12] = 0.9941756912993762

# row 11 - Spiral-90deg-10rot-5space
$data[10,0] = 9
$data[10,1] = 'Spiral-90deg-10rot-5space'
$data[10,2] = 0.9991865085824876
$data[10,3] = 0.9533782281102589
$data[10,4] = 1.006861331950889
$data[10,5] = 0.9991865085824876
$data[10,6] = 0.9736533217704941
$data[10,7] = 1.027087365986937
$data[10,8] = 1.006211609223991
$data[10,9] = 0.9533782281102589
$data[10,10] = 0.9801197800305739
$data[10,11] = 0.9896531443065307
$data[10,12] = 0.9943963942708428

# row 12 - Spiral-90deg-15rot-5space
$data[11,0] = 10
$data[11,1] = 'Spiral-90deg-15rot-5space'
$data[11,2] = 0.9990954443368405
$data[11,3] = 0.95381197143259
$data[11,4] = 1.006766363654171
$data[11,5] = 0.9990954443368405
$data[11,6] = 0.9738470336726458
$data[11,7] = 1.026840782890429
$data[11,8] = 1.006137806986086
$data[11,9] = 0.95381197143259
$data[11,10] = 0.9802891675433802
$data[11,11] = 0.9896923059401104
$data[11,12] = 0.9944165671621269

# row 13 - Spiral-90deg-10rot-3space
$data[12,0] = 11
$data[12,1] = 'Spiral-90deg-10rot-3space'
$data[12,2] = 0.9991600302691497
$data[12,3] = 0.9534933989462824
$data[12,4] = 1.006834044073103
$data[12,5] = 0.9991600302691497
$data[12,6] = 0.9736901735243946
$data[12,7] = 1.026986321064178
$data[12,8] = 1.006213451167253
$data[12,9] = 0.9534933989462824
$data[12,10] = 0.9801637215096924
$data[12,11] = 0.9896618758894211
$data[12,12] = 0.9943962365073933

# row 14 - NoRotation-tilt60deg
$data[13,0] = 12
$data[13,1] = 'NoRotation-tilt60deg'
$data[13,2] = 0.9827160000000005
$data[13,3] = 1.019967999999999
$data[13,4] = 0.9918840000000002
$data[13,5] = 0.9827160000000005
$data[13,6] = 1.012903999999999
$data[13,7] = 0.9958559999999994
$data[13,8] = 0.9860440000000003
$data[13,9] = 1.019967999999999
$data[13,10] = 1.005925999999999
$data[13,11] = 0.9943209999999999
$data[13,12] = 0.9982286666666664

# row 15 - Rotation-NoTilt
$data[14,0] = 13
$data[14,1] = 'Rotation-NoTilt'
$data[14,2] = 0.98
$data[14,3] = 1.03
$data[14,4] = 0.99
$data[14,5] = 0.98
$data[14,6] = 1.02
$data[14,7] = 1
$data[14,8] = 0.98
$data[14,9] = 1.03
$data[14,10] = 1.01
$data[14,11] = 0.995
$data[14,12] = 1

# row 16 - Rotation-60detTilt
$data[15,0] = 14
$data[15,1] = 'Rotation-60detTilt'
$data[15,2] = 0.9862614800384024
$data[15,3] = 1.015557374156799
$data[15,4] = 0.9919629619200008
$data[15,5] = 0.9862614800384024
$data[15,6] = 1.008705201356798
$data[15,7] = 0.9968439814143992
$data[15,8] = 0.9869389862912015
$data[15,9] = 1.015557374156799
$data[15,10] = 1.0037601680384
$data[15,11] = 0.9950108240384012
$data[15,12] = 0.997711664196267

# row 17 - HexGrid-90degTilt5degRes
$data[16,0] = 15
$data[16,1] = 'HexGrid-90degTilt5degRes'
$data[16,2] = 0.994884828880118
$data[16,3] = 0.9946628661814049
$data[16,4] = 0.9950140723697587
$data[16,5] = 0.994884828880118
$data[16,6] = 0.9949225688599869
$data[16,7] = 0.9948529628041923
$data[16,8] = 0.9951108940281472
$data[16,9] = 0.9946628661814049
$data[16,10] = 0.9948384692755818
$data[16,11] = 0.9948616490778499
$data[16,12] = 0.994908032187268

# row 18 - HexGrid-90degTilt22p5degRes
$data[17,0] = 16
$data[17,1] = 'HexGrid-90degTilt22p5degRes'
$data[17,2] = 0.9948960510810921
$data[17,3] = 0.9928922693160657
$data[17,4] = 0.9955563243288749
$data[17,5] = 0.9948960510810921
$data[17,6] = 0.9942844847104138
$data[17,7] = 0.9948279706641349
$data[17,8] = 0.9967402599771671
$data[17,9] = 0.9928922693160657
$data[17,10] = 0.9942242968224703
$data[17,11] = 0.9945601739517812
$data[17,12] = 0.9948662266796248

# row 19 - HexGrid-60degTilt5degRes
$data[18,0] = 17
$data[18,1] = 'HexGrid-60degTilt5degRes'
$data[18,2] = 0.9958306887590581
$data[18,3] = 0.9916499695233784
$data[18,4] = 0.9955543055001368
$data[18,5] = 0.9958306887590581
$data[18,6] = 0.9924642169604299
$data[18,7] = 0.9968278875797733
$data[18,8] = 0.9962320654084875
$data[18,9] = 0.9916499695233784
$data[18,10] = 0.9936021375117576
$data[18,11] = 0.9947164131354078
$data[18,12] = 0.9947598556218772

$ws.Range("A1:M19").Value = $data

# The three newly appended rows (17-19) need the same bold/bordered/centered
# label style already used by the other rows in column A; grab it from the
# last pre-existing styled row instead of hard-coding the format.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
